$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("1er Parcial")
$ws1.Range("I13").Value = 9.199999999999999

$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Range("E2").Value = 38
$ws2.Range("F2").Value = 1
$ws2.Range("G2").Value = 97.40000000000001
$ws2.Range("H2").Value = 2.6
$ws2.Range("I2").Value = 8.5
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0
$ws2.Range("E3").Value = 38
$ws2.Range("F3").Value = 1
$ws2.Range("G3").Value = 97.40000000000001
$ws2.Range("H3").Value = 2.6
$ws2.Range("I3").Value = 8.4
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0
$ws2.Range("E4").Value = 33
$ws2.Range("F4").Value = 5
$ws2.Range("G4").Value = 86.8
$ws2.Range("H4").Value = 13.2
$ws2.Range("I4").Value = 7.6
$ws2.Range("J4").Value = 0
$ws2.Range("K4").Value = 0
$ws2.Range("E5").Value = 109
$ws2.Range("F5").Value = 7
$ws2.Range("G5").Value = 94
$ws2.Range("H5").Value = 6
$ws2.Range("I5").Value = 8.199999999999999
$ws2.Range("J5").Value = 0
$ws2.Range("K5").Value = 0
$ws2.Range("E6").Value = 26
$ws2.Range("F6").Value = 2
$ws2.Range("G6").Value = 92.90000000000001
$ws2.Range("H6").Value = 7.1
$ws2.Range("I6").Value = 7.9
$ws2.Range("J6").Value = 0
$ws2.Range("K6").Value = 0
$ws2.Range("E7").Value = 26
$ws2.Range("F7").Value = 2
$ws2.Range("G7").Value = 92.90000000000001
$ws2.Range("H7").Value = 7.1
$ws2.Range("I7").Value = 8.699999999999999
$ws2.Range("J7").Value = 0
$ws2.Range("K7").Value = 0
$ws2.Range("E8").Value = 52
$ws2.Range("F8").Value = 4
$ws2.Range("G8").Value = 92.90000000000001
$ws2.Range("H8").Value = 7.1
$ws2.Range("I8").Value = 8.300000000000001
$ws2.Range("J8").Value = 0
$ws2.Range("K8").Value = 0
$ws2.Range("E9").Value = 31
$ws2.Range("F9").Value = 7
$ws2.Range("G9").Value = 81.59999999999999
$ws2.Range("H9").Value = 18.4
$ws2.Range("I9").Value = 8
$ws2.Range("J9").Value = 0
$ws2.Range("K9").Value = 0
$ws2.Range("E10").Value = 34
$ws2.Range("F10").Value = 0
$ws2.Range("G10").Value = 100
$ws2.Range("H10").Value = 0
$ws2.Range("I10").Value = 9
$ws2.Range("J10").Value = 0
$ws2.Range("K10").Value = 0
$ws2.Range("E11").Value = 34
$ws2.Range("F11").Value = 0
$ws2.Range("G11").Value = 100
$ws2.Range("H11").Value = 0
$ws2.Range("I11").Value = 9.300000000000001
$ws2.Range("J11").Value = 0
$ws2.Range("K11").Value = 0
$ws2.Range("E12").Value = 38
$ws2.Range("F12").Value = 0
$ws2.Range("G12").Value = 100
$ws2.Range("H12").Value = 0
$ws2.Range("I12").Value = 9.199999999999999
$ws2.Range("J12").Value = 0
$ws2.Range("K12").Value = 0
$ws2.Range("E13").Value = 38
$ws2.Range("F13").Value = 0
$ws2.Range("G13").Value = 100
$ws2.Range("H13").Value = 0
$ws2.Range("I13").Value = 8.9
$ws2.Range("J13").Value = 0
$ws2.Range("K13").Value = 0
$ws2.Range("E14").Value = 38
$ws2.Range("F14").Value = 0
$ws2.Range("G14").Value = 100
$ws2.Range("H14").Value = 0
$ws2.Range("I14").Value = 9.199999999999999
$ws2.Range("J14").Value = 0
$ws2.Range("K14").Value = 0
$ws2.Range("E15").Value = 213
$ws2.Range("F15").Value = 7
$ws2.Range("G15").Value = 96.8
$ws2.Range("H15").Value = 3.2
$ws2.Range("I15").Value = 8.9
$ws2.Range("J15").Value = 0
$ws2.Range("K15").Value = 0
$ws2.Range("E16").Value = 24
$ws2.Range("F16").Value = 4
$ws2.Range("G16").Value = 85.7
$ws2.Range("H16").Value = 14.3
$ws2.Range("I16").Value = 7.5
$ws2.Range("J16").Value = 0
$ws2.Range("K16").Value = 0
$ws2.Range("E17").Value = 25
$ws2.Range("F17").Value = 0
$ws2.Range("G17").Value = 100
$ws2.Range("H17").Value = 0
$ws2.Range("I17").Value = 7.7
$ws2.Range("J17").Value = 0
$ws2.Range("K17").Value = 0
$ws2.Range("E18").Value = 49
$ws2.Range("F18").Value = 4
$ws2.Range("G18").Value = 92.5
$ws2.Range("H18").Value = 7.5
$ws2.Range("I18").Value = 7.6
$ws2.Range("J18").Value = 0
$ws2.Range("K18").Value = 0
$ws2.Range("E19").Value = 25
$ws2.Range("F19").Value = 0
$ws2.Range("G19").Value = 100
$ws2.Range("H19").Value = 0
$ws2.Range("I19").Value = 8.4
$ws2.Range("J19").Value = 0
$ws2.Range("K19").Value = 0
$ws2.Range("E20").Value = 25
$ws2.Range("F20").Value = 0
$ws2.Range("G20").Value = 100
$ws2.Range("H20").Value = 0
$ws2.Range("I20").Value = 8.4
$ws2.Range("J20").Value = 0
$ws2.Range("K20").Value = 0
$ws2.Range("E21").Value = 27
$ws2.Range("F21").Value = 1
$ws2.Range("G21").Value = 96.40000000000001
$ws2.Range("H21").Value = 3.6
$ws2.Range("I21").Value = 8.199999999999999
$ws2.Range("J21").Value = 0
$ws2.Range("K21").Value = 0
$ws2.Range("E22").Value = 77
$ws2.Range("F22").Value = 1
$ws2.Range("G22").Value = 98.7
$ws2.Range("H22").Value = 1.3
$ws2.Range("I22").Value = 8.300000000000001
$ws2.Range("J22").Value = 0
$ws2.Range("K22").Value = 0
$ws2.Range("E23").Value = 34
$ws2.Range("F23").Value = 0
$ws2.Range("G23").Value = 100
$ws2.Range("H23").Value = 0
$ws2.Range("I23").Value = 9.5
$ws2.Range("J23").Value = 0
$ws2.Range("K23").Value = 0
$ws2.Range("E24").Value = 40
$ws2.Range("F24").Value = 1
$ws2.Range("G24").Value = 97.59999999999999
$ws2.Range("H24").Value = 2.4
$ws2.Range("I24").Value = 8.800000000000001
$ws2.Range("J24").Value = 0
$ws2.Range("K24").Value = 0
$ws2.Range("E25").Value = 40
$ws2.Range("F25").Value = 1
$ws2.Range("G25").Value = 97.59999999999999
$ws2.Range("H25").Value = 2.4
$ws2.Range("I25").Value = 8.800000000000001
$ws2.Range("J25").Value = 0
$ws2.Range("K25").Value = 0
$ws2.Range("E26").Value = 36
$ws2.Range("F26").Value = 0
$ws2.Range("G26").Value = 100
$ws2.Range("H26").Value = 0
$ws2.Range("I26").Value = 9
$ws2.Range("J26").Value = 0
$ws2.Range("K26").Value = 0
$ws2.Range("E27").Value = 36
$ws2.Range("F27").Value = 0
$ws2.Range("G27").Value = 100
$ws2.Range("H27").Value = 0
$ws2.Range("I27").Value = 9
$ws2.Range("J27").Value = 0
$ws2.Range("K27").Value = 0
$ws2.Range("E28").Value = 186
$ws2.Range("F28").Value = 2
$ws2.Range("G28").Value = 98.90000000000001
$ws2.Range("H28").Value = 1.1
$ws2.Range("I28").Value = 9
$ws2.Range("J28").Value = 0
$ws2.Range("K28").Value = 0
$ws2.Range("E29").Value = 686
$ws2.Range("F29").Value = 25
$ws2.Range("G29").Value = 96.5
$ws2.Range("H29").Value = 3.5
$ws2.Range("I29").Value = 8.6
$ws2.Range("J29").Value = 0
$ws2.Range("K29").Value = 0

$ws3 = $wb.Worksheets.Item("Final")
$ws3.Range("E2").Value = 38
$ws3.Range("F2").Value = 1
$ws3.Range("G2").Value = 97.40000000000001
$ws3.Range("H2").Value = 2.6
$ws3.Range("I2").Value = 8.1
$ws3.Range("I3").Value = 8.800000000000001
$ws3.Range("E4").Value = 33
$ws3.Range("F4").Value = 5
$ws3.Range("G4").Value = 86.8
$ws3.Range("H4").Value = 13.2
$ws3.Range("I4").Value = 8.199999999999999
$ws3.Range("E5").Value = 109
$ws3.Range("F5").Value = 7
$ws3.Range("G5").Value = 94
$ws3.Range("H5").Value = 6
$ws3.Range("I5").Value = 8.4
$ws3.Range("E6").Value = 26
$ws3.Range("F6").Value = 2
$ws3.Range("G6").Value = 92.90000000000001
$ws3.Range("H6").Value = 7.1
$ws3.Range("I6").Value = 7.6
$ws3.Range("E7").Value = 26
$ws3.Range("F7").Value = 2
$ws3.Range("G7").Value = 92.90000000000001
$ws3.Range("H7").Value = 7.1
$ws3.Range("I7").Value = 8.6
$ws3.Range("E8").Value = 52
$ws3.Range("F8").Value = 4
$ws3.Range("G8").Value = 92.90000000000001
$ws3.Range("H8").Value = 7.1
$ws3.Range("I8").Value = 8.1
$ws3.Range("E9").Value = 31
$ws3.Range("F9").Value = 7
$ws3.Range("G9").Value = 81.59999999999999
$ws3.Range("H9").Value = 18.4
$ws3.Range("I9").Value = 7.9
$ws3.Range("I10").Value = 8.9
$ws3.Range("I11").Value = 9.4
$ws3.Range("I12").Value = 9.199999999999999
$ws3.Range("I13").Value = 9.199999999999999
$ws3.Range("I14").Value = 9.4
$ws3.Range("E15").Value = 213
$ws3.Range("F15").Value = 7
$ws3.Range("G15").Value = 96.8
$ws3.Range("H15").Value = 3.2
$ws3.Range("I15").Value = 9
$ws3.Range("I16").Value = 7.3
$ws3.Range("I17").Value = 7.8
$ws3.Range("I18").Value = 7.6
$ws3.Range("E19").Value = 25
$ws3.Range("F19").Value = 0
$ws3.Range("G19").Value = 100
$ws3.Range("H19").Value = 0
$ws3.Range("I19").Value = 7.4
$ws3.Range("E20").Value = 25
$ws3.Range("F20").Value = 0
$ws3.Range("G20").Value = 100
$ws3.Range("H20").Value = 0
$ws3.Range("I20").Value = 8.1
$ws3.Range("E21").Value = 27
$ws3.Range("F21").Value = 1
$ws3.Range("G21").Value = 96.40000000000001
$ws3.Range("H21").Value = 3.6
$ws3.Range("I21").Value = 8.4
$ws3.Range("E22").Value = 77
$ws3.Range("F22").Value = 1
$ws3.Range("G22").Value = 98.7
$ws3.Range("H22").Value = 1.3
$ws3.Range("I22").Value = 8
$ws3.Range("I23").Value = 9.199999999999999
$ws3.Range("I24").Value = 9
$ws3.Range("I25").Value = 9
$ws3.Range("I26").Value = 9.300000000000001
$ws3.Range("I27").Value = 9.199999999999999
$ws3.Range("I28").Value = 9.1
$ws3.Range("E29").Value = 686
$ws3.Range("F29").Value = 25
$ws3.Range("G29").Value = 96.5
$ws3.Range("H29").Value = 3.5
$ws3.Range("I29").Value = 8.6
